# Apply the weekly hortaliza/fruta update:
# Insert a new row at position 70 (old row 70 data shifts down to row 71),
# then populate the new row 70 with the latest weekly price record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift row 70 (and everything below) down by one row, duplicating formatting.
$ws.Rows.Item(70).Insert()

# Populate the newly inserted row 70 with the new weekly data.
$ws.Range("A70").Value = 1
$ws.Range("B70").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C70").Value = "Arica y Parinacota"
$ws.Range("D70").Value = 44747
$ws.Range("E70").Value = 15
$ws.Range("F70").Value = 100112012
$ws.Range("G70").Value = "Espinaca"
$ws.Range("H70").Value = "Sin especificar"
$ws.Range("I70").Value = "Primera"
$ws.Range("J70").Value = 200
$ws.Range("K70").Value = 2000
$ws.Range("L70").Value = 2500
$ws.Range("M70").Value = 2250
$ws.Range("N70").Value = "$/atado 2,5 a 3 kilos"
$ws.Range("O70").Value = "Región de Arica y Parinacota"
$ws.Range("P70").Value = 750
$ws.Range("Q70").Value = 3
$ws.Range("R70").Value = "Hortaliza"

$wb.Save()
